$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.417.81'
$ws.Range('D3').Value = '1.867.94'
$ws.Range('E3').Value = '  -0.59%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').Value = "'" + '243.73'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.51%  '
$ws.Range('B6').Value = 'XRP'
$ws.Range('C6').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D6').Value = "'" + '0.7069'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.60%  '
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').Value = "'" + '0.07861'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.08%  '
$ws.Range('D9').Value = "'" + '0.3132'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.98%  '
$ws.Range('D10').Value = "'" + '24.50'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.98%  '
$ws.Range('D11').Value = "'" + '0.07969'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.07%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.887.68'
$ws.Range('E12').Value = '  +0.41%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = "'" + '5.209'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.84%  '
$ws.Range('D14').Value = "'" + '93.35'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E15').Value = '  -2.12%  '
$ws.Range('D16').Value = "'" + '6.518'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.61%  '
$ws.Range('D17').Value = "'" + '0.000008367'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.03%  '
$ws.Range('D18').Value = '29.428.94'
$ws.Range('E18').Value = '  +0.27%  '
$ws.Range('D19').Value = "'" + '252.55'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.62%  '
$ws.Range('D20').Value = '2.124.53'
$ws.Range('E20').Value = '  -0.34%  '
$ws.Range('D21').Value = "'" + '13.12'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.37%  '
$ws.Range('E22').Value = '  -0.11%  '
$ws.Range('D23').Value = "'" + '7.630'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.25%  '
$ws.Range('E24').Value = '  -0.30%  '
$ws.Range('D25').Value = "'" + '0.1552'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.63%  '
$ws.Range('D26').Value = "'" + '9.002'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Value = "'" + '161.23'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.02%  '
$ws.Range('D28').Value = "'" + '18.70'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Value = "'" + '1.500'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.51%  '
$ws.Range('E30').Value = '  -2.22%  '
$ws.Range('D31').Value = "'" + '4.258'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.60%  '
$ws.Range('E32').Value = '  +1.37%  '
$ws.Range('D33').Value = "'" + '0.05318'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.24%  '
$ws.Range('D34').Value = "'" + '1.892'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.33%  '
$ws.Range('D35').Value = "'" + '0.7486'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.84%  '
$ws.Range('E36').Value = '  -1.16%  '
$ws.Range('D37').Value = "'" + '2.711'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.99%  '
$ws.Range('D38').Value = "'" + '0.01889'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.17%  '
$ws.Range('D39').Value = '1.274.63'
$ws.Range('E39').Value = '  +0.96%  '
$ws.Range('D40').Value = "'" + '2.748'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Value = "'" + '0.8924'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.58%  '
$ws.Range('D42').Value = "'" + '6.056'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -6.95%  '
$ws.Range('D43').Value = "'" + '109.00'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.78%  '
$ws.Range('D44').Value = "'" + '71.23'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.19%  '
$ws.Range('E45').Value = '  -0.11%  '
$ws.Range('E46').Value = '  -4.06%  '
$ws.Range('D47').Value = '2.019.26'
$ws.Range('E47').Value = '  -0.11%  '
$ws.Range('E48').Value = '  -0.40%  '
$ws.Range('D49').Value = "'" + '9.552'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.99%  '
$ws.Range('D50').Value = "'" + '0.5179'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.99%  '
$ws.Range('D51').Value = "'" + '0.4305'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.72%  '
